$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 45.76760875408854
$ws.Range("B3").Value = 32.08920247651583
$ws.Range("B4").Value = 12.59201073713469
$ws.Range("B5").Value = 9.551178032260943
